# Add a "blood_type list" sheet (enum source) after the existing data sheet,
# populate it with the four blood-type values, and wire a List data
# validation on the "blood_type" column (B) of the main sheet back to it.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1. Create the lookup sheet, placed after the main sheet -------------
$listSheet = $wb.Worksheets.Add($null, $ws1)
$listSheet.Name = "blood_type list"

$listSheet.Range("A1").Value = "A"
$listSheet.Range("A2").Value = "B"
$listSheet.Range("A3").Value = "AB"
$listSheet.Range("A4").Value = "O"

# --- 2. Re-create the data validations on the main sheet so the new list
#        validation for column B lands between the A and C validations,
#        matching the column order of the sheet. ---------------------------

function Set-DecimalValidation($rangeAddr) {
    $v = $ws1.Range($rangeAddr).Validation
    $v.Delete()
    $v.Add(2, 1, 1, "-1e+307", "1e+307")
    $v.ErrorTitle = "Not a number"
    $v.ErrorMessage = "The values in this column must be numbers."
}

Set-DecimalValidation("A2:A1048576")

$bloodValidation = $ws1.Range("B2:B1048576").Validation
$bloodValidation.Add(3, 1, 1, "='blood_type list'!`$A`$1:`$A`$4")
$bloodValidation.ErrorTitle = "Value must come from list"
$bloodValidation.ErrorMessage = "Value must be one of: A / B / AB / O."

Set-DecimalValidation("C2:C1048576")
Set-DecimalValidation("F2:F1048576")
Set-DecimalValidation("K2:K1048576")

# --- 3. Leave the original sheet as the active/selected one, as before ---
$ws1.Activate()
